$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cell values (B2:J35) per diff
$ws.Range("B2").Value = 0.008403361344537815
$ws.Range("D2").Value = 0
$ws.Range("F2").Value = 0.01149425287356322
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0.02478667208451852
$ws.Range("I2").Value = 0.01186943620178042
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0.1417092768444125
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0.02155172413793103
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0.03657049979683057
$ws.Range("I3").Value = 0.02670623145400593
$ws.Range("J3").Value = 0.003073140749846343
$ws.Range("B4").Value = 0.008403361344537815
$ws.Range("C4").Value = 0.004626247869491115
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0.007203667321545513
$ws.Range("F4").Value = 0.008620689655172414
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0.006095083299471759
$ws.Range("I4").Value = 0.01483679525222552
$ws.Range("J4").Value = 0.0313460356484327
$ws.Range("B5").Value = 0.3361344537815129
$ws.Range("C5").Value = 0.001217433649866082
$ws.Range("D5").Value = 1
$ws.Range("F5").Value = 0.2514367816091963
$ws.Range("G5").Value = 0.625
$ws.Range("H5").Value = 0.01137748882568062
$ws.Range("I5").Value = 0.3738872403560825
$ws.Range("J5").Value = 0
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0.01266130995860726
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0.01724137931034483
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0.01625355546525804
$ws.Range("I6").Value = 0.01483679525222552
$ws.Range("C7").Value = 0.002921840759678598
$ws.Range("E7").Value = 0.04584151931892596
$ws.Range("H7").Value = 0.0008126777732629012
$ws.Range("J7").Value = 0.02888752304855562
$ws.Range("B8").Value = 0.0588235294117647
$ws.Range("C8").Value = 0.03116630143657162
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0.1473477406679767
$ws.Range("F8").Value = 0.06609195402298845
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0.03860219422998781
$ws.Range("I8").Value = 0.06231454005934716
$ws.Range("J8").Value = 0.1143208358942835
$ws.Range("C9").Value = 0.007061115169223281
$ws.Range("E9").Value = 0.09692206941715807
$ws.Range("H9").Value = 0.002844372206420154
$ws.Range("J9").Value = 0.07129686539643494
$ws.Range("B10").Value = 0.008403361344537815
$ws.Range("C10").Value = 0.0009739469198928658
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0.07531106745252134
$ws.Range("F10").Value = 0.004310344827586207
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0.002438033319788704
$ws.Range("J10").Value = 0.01413644744929318
$ws.Range("E11").Value = 0.01375245579567781
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0.003073140749846343
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0.0007304601899196494
$ws.Range("E12").Value = 0.02750491159135559
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0.01352181929932391
$ws.Range("B13").Value = 0.06722689075630252
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("F13").Value = 0.03160919540229883
$ws.Range("G13").Value = 0.125
$ws.Range("H13").Value = 0.007720438845997562
$ws.Range("I13").Value = 0.02670623145400593
$ws.Range("J13").Value = 0
$ws.Range("B14").Value = 0.008403361344537815
$ws.Range("D14").Value = 0
$ws.Range("F14").Value = 0.1020114942528735
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0.002438033319788704
$ws.Range("I14").Value = 0.02967359050445104
$ws.Range("B16").Value = 0.0588235294117647
$ws.Range("C16").Value = 0.04699293888483059
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0.04846103470857887
$ws.Range("F16").Value = 0.08333333333333326
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0.08573750507923616
$ws.Range("I16").Value = 0.07418397626112759
$ws.Range("J16").Value = 0.1020282728948982
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = 0.0131482834185537
$ws.Range("D17").Value = 0
$ws.Range("F17").Value = 0.004310344827586207
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0.01747257212515239
$ws.Range("E18").Value = 0.01113294040602489
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0.0006146281499692685
$ws.Range("B19").Value = 0.1260504201680672
$ws.Range("D19").Value = 0
$ws.Range("F19").Value = 0.03304597701149423
$ws.Range("G19").Value = 0.125
$ws.Range("I19").Value = 0.09198813056379825
$ws.Range("J19").Value = 0
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 0.03944485025566093
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0.001964636542239686
$ws.Range("F20").Value = 0.02729885057471263
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 0.1036164160910202
$ws.Range("I20").Value = 0.01483679525222552
$ws.Range("J20").Value = 0.01044867854947757
$ws.Range("B21").Value = 0.04201680672268907
$ws.Range("C21").Value = 0.04358412466520558
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0.01702685003274395
$ws.Range("F21").Value = 0.03304597701149423
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0.05485574969524572
$ws.Range("I21").Value = 0.03560830860534124
$ws.Range("J21").Value = 0.08113091579594318
$ws.Range("B23").Value = 0.008403361344537815
$ws.Range("C23").Value = 0.1190650109569038
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0.01702685003274395
$ws.Range("F23").Value = 0.05459770114942524
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0.05241771637545704
$ws.Range("I23").Value = 0.03264094955489614
$ws.Range("J23").Value = 0.1032575291948367
$ws.Range("E24").Value = 0.01964636542239686
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0.002458512599877074
$ws.Range("B25").Value = 0
$ws.Range("C25").Value = 0.006817628439250064
$ws.Range("E25").Value = 0.003274394237066143
$ws.Range("H25").Value = 0.002438033319788704
$ws.Range("J25").Value = 0.05285802089735698
$ws.Range("E26").Value = 0.0006548788474132286
$ws.Range("H26").Value = 0
$ws.Range("E27").Value = 0.03536345776031432
$ws.Range("H27").Value = 0
$ws.Range("B28").Value = 0
$ws.Range("C28").Value = 0.01241782322863405
$ws.Range("E28").Value = 0
$ws.Range("H28").Value = 0.0004063388866314506
$ws.Range("J28").Value = 0.03749231714812535
$ws.Range("B29").Value = 0
$ws.Range("C29").Value = 0.001947893839785732
$ws.Range("E29").Value = 0.01964636542239686
$ws.Range("H29").Value = 0.001219016659894352
$ws.Range("J29").Value = 0.01413644744929318
$ws.Range("B30").Value = 0
$ws.Range("C30").Value = 0.002921840759678598
$ws.Range("E30").Value = 0.04780615586116564
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0.002844372206420154
$ws.Range("J30").Value = 0.009834050399508297
$ws.Range("B31").Value = 0
$ws.Range("C31").Value = 0.002191380569758948
$ws.Range("E31").Value = 0.04256712508185982
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0.02397049784880147
$ws.Range("B32").Value = 0.02521008403361345
$ws.Range("C32").Value = 0.03603603603603592
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0.07334643091028163
$ws.Range("F32").Value = 0.02442528735632183
$ws.Range("G32").Value = 0.0625
$ws.Range("H32").Value = 0.02234863876472981
$ws.Range("I32").Value = 0.01483679525222552
$ws.Range("J32").Value = 0.1567301782421628
$ws.Range("B33").Value = 0
$ws.Range("C33").Value = 0.007304601899196498
$ws.Range("E33").Value = 0.1656843483955468
$ws.Range("F33").Value = 0.004310344827586207
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 0.005688744412840309
$ws.Range("J33").Value = 0.05162876459741845
$ws.Range("B34").Value = 0
$ws.Range("D34").Value = 0
$ws.Range("F34").Value = 0.01149425287356322
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0.01381552214546933
$ws.Range("I34").Value = 0.002967359050445104
$ws.Range("E35").Value = 0.03994760969220691
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0.0006146281499692685

# Remove now-unused "Joint regime area" rows (36-40); dimension shrinks to A1:J35 automatically
$ws.Range("A36:A40").EntireRow.Delete()
